# Scheduled data refresh: updates market-price derived columns
# (currentAveragePrice/NQ/HQ, LevePrice NQ/HQ, LeveProfit NQ/HQ)
# for a handful of leve rows across the per-job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 3979.2083
$ws.Cells.Item(64, 9).Value = 3834.7827
$ws.Cells.Item(64, 10).Value = 4112.08
$ws.Cells.Item(64, 11).Value = 3834.7827
$ws.Cells.Item(64, 12).Value = 4112.08
$ws.Cells.Item(64, 13).Value = -3586.7827
$ws.Cells.Item(64, 14).Value = -4608.08

$ws.Cells.Item(67, 8).Value = 3979.2083
$ws.Cells.Item(67, 9).Value = 3834.7827
$ws.Cells.Item(67, 10).Value = 4112.08
$ws.Cells.Item(67, 11).Value = 3834.7827
$ws.Cells.Item(67, 12).Value = 4112.08
$ws.Cells.Item(67, 13).Value = -2976.7827
$ws.Cells.Item(67, 14).Value = -5828.08

$ws.Cells.Item(70, 8).Value = 5130213
$ws.Cells.Item(70, 9).Value = 1525.5
$ws.Cells.Item(70, 10).Value = 7409629.5
$ws.Cells.Item(70, 11).Value = 4576.5
$ws.Cells.Item(70, 12).Value = 22228888.5
$ws.Cells.Item(70, 13).Value = -4306.5
$ws.Cells.Item(70, 14).Value = -22229428.5

$ws.Cells.Item(73, 8).Value = 5130213
$ws.Cells.Item(73, 9).Value = 1525.5
$ws.Cells.Item(73, 10).Value = 7409629.5
$ws.Cells.Item(73, 11).Value = 4576.5
$ws.Cells.Item(73, 12).Value = 22228888.5
$ws.Cells.Item(73, 13).Value = -3640.5
$ws.Cells.Item(73, 14).Value = -22230760.5

$ws.Cells.Item(92, 8).Value = 3473.9
$ws.Cells.Item(92, 9).Value = 864
$ws.Cells.Item(92, 10).Value = 7388.75
$ws.Cells.Item(92, 11).Value = 864
$ws.Cells.Item(92, 12).Value = 7388.75
$ws.Cells.Item(92, 13).Value = 384
$ws.Cells.Item(92, 14).Value = -9884.75

$ws.Cells.Item(129, 8).Value = 979.45
$ws.Cells.Item(129, 9).Value = 379.58334
$ws.Cells.Item(129, 10).Value = 1061.25
$ws.Cells.Item(129, 11).Value = 1138.75002
$ws.Cells.Item(129, 12).Value = 3183.75
$ws.Cells.Item(129, 13).Value = 3861.24998
$ws.Cells.Item(129, 14).Value = -13183.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(102, 8).Value = 3232.6365
$ws.Cells.Item(102, 9).Value = 3232.6365
$ws.Cells.Item(102, 10).Value = 0
$ws.Cells.Item(102, 11).Value = 3232.6365
$ws.Cells.Item(102, 12).Value = 0
$ws.Cells.Item(102, 13).Value = -1610.6365
$ws.Cells.Item(102, 14).Value = $null

$ws.Cells.Item(122, 8).Value = 1849.5555
$ws.Cells.Item(122, 9).Value = 1852.0667
$ws.Cells.Item(122, 10).Value = 1837
$ws.Cells.Item(122, 11).Value = 5556.2001
$ws.Cells.Item(122, 12).Value = 5511
$ws.Cells.Item(122, 13).Value = -3106.2001
$ws.Cells.Item(122, 14).Value = -10411

$ws.Cells.Item(123, 8).Value = 34984.5
$ws.Cells.Item(123, 10).Value = 34984.5
$ws.Cells.Item(123, 12).Value = 34984.5
$ws.Cells.Item(123, 14).Value = -44784.5

$ws.Cells.Item(132, 8).Value = 1998.2858
$ws.Cells.Item(132, 9).Value = 1432.4706
$ws.Cells.Item(132, 11).Value = 4297.4118
$ws.Cells.Item(132, 13).Value = -1767.4118

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 2772.41
$ws.Cells.Item(105, 9).Value = 1785.5714
$ws.Cells.Item(105, 10).Value = 2846.6882
$ws.Cells.Item(105, 11).Value = 1785.5714
$ws.Cells.Item(105, 12).Value = 2846.6882
$ws.Cells.Item(105, 13).Value = -38.57140000000004
$ws.Cells.Item(105, 14).Value = -6340.688200000001

$ws.Cells.Item(107, 8).Value = 3420.7856
$ws.Cells.Item(107, 9).Value = 3301.3333
$ws.Cells.Item(107, 10).Value = 3635.8
$ws.Cells.Item(107, 11).Value = 3301.3333
$ws.Cells.Item(107, 12).Value = 3635.8
$ws.Cells.Item(107, 13).Value = -1381.3333
$ws.Cells.Item(107, 14).Value = -7475.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(107, 8).Value = 498.23914
$ws.Cells.Item(107, 9).Value = 457.6
$ws.Cells.Item(107, 10).Value = 574.4375
$ws.Cells.Item(107, 11).Value = 457.6
$ws.Cells.Item(107, 12).Value = 574.4375
$ws.Cells.Item(107, 13).Value = 1462.4
$ws.Cells.Item(107, 14).Value = -4414.4375

$ws.Cells.Item(141, 8).Value = 27902.375
$ws.Cells.Item(141, 9).Value = 21000
$ws.Cells.Item(141, 10).Value = 29495.23
$ws.Cells.Item(141, 11).Value = 21000
$ws.Cells.Item(141, 12).Value = 29495.23
$ws.Cells.Item(141, 13).Value = -15820
$ws.Cells.Item(141, 14).Value = -39855.23

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(129, 8).Value = 6762.75
$ws.Cells.Item(129, 9).Value = 2294.875
$ws.Cells.Item(129, 10).Value = 9741.333000000001
$ws.Cells.Item(129, 11).Value = 6884.625
$ws.Cells.Item(129, 12).Value = 29223.999
$ws.Cells.Item(129, 13).Value = -1884.625
$ws.Cells.Item(129, 14).Value = -39223.999

$ws.Cells.Item(133, 8).Value = 2913.889
$ws.Cells.Item(133, 9).Value = 1693.1666
$ws.Cells.Item(133, 10).Value = 5355.3335
$ws.Cells.Item(133, 11).Value = 5079.4998
$ws.Cells.Item(133, 12).Value = 16066.0005
$ws.Cells.Item(133, 13).Value = -19.4997999999996
$ws.Cells.Item(133, 14).Value = -26186.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(93, 8).Value = 82450
$ws.Cells.Item(93, 10).Value = 82450
$ws.Cells.Item(93, 12).Value = 82450
$ws.Cells.Item(93, 14).Value = -86194

$ws.Cells.Item(102, 8).Value = 1611.1786
$ws.Cells.Item(102, 9).Value = 1569.0416
$ws.Cells.Item(102, 10).Value = 1864
$ws.Cells.Item(102, 11).Value = 1569.0416
$ws.Cells.Item(102, 12).Value = 1864
$ws.Cells.Item(102, 13).Value = 52.95839999999998
$ws.Cells.Item(102, 14).Value = -5108

$ws.Cells.Item(132, 8).Value = 3896
$ws.Cells.Item(132, 9).Value = 3235
$ws.Cells.Item(132, 10).Value = 4997.6665
$ws.Cells.Item(132, 11).Value = 9705
$ws.Cells.Item(132, 12).Value = 14992.9995
$ws.Cells.Item(132, 13).Value = -7175
$ws.Cells.Item(132, 14).Value = -20052.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2589.5789
$ws.Cells.Item(7, 9).Value = 2587.4375
$ws.Cells.Item(7, 10).Value = 2601
$ws.Cells.Item(7, 11).Value = 2587.4375
$ws.Cells.Item(7, 12).Value = 2601
$ws.Cells.Item(7, 13).Value = -2475.4375
$ws.Cells.Item(7, 14).Value = -2825

$ws.Cells.Item(40, 8).Value = 5275.5625
$ws.Cells.Item(40, 9).Value = 5825.5
$ws.Cells.Item(40, 10).Value = 4725.625
$ws.Cells.Item(40, 11).Value = 5825.5
$ws.Cells.Item(40, 12).Value = 4725.625
$ws.Cells.Item(40, 13).Value = -5689.5
$ws.Cells.Item(40, 14).Value = -4997.625

$ws.Cells.Item(87, 8).Value = 42250
$ws.Cells.Item(87, 10).Value = 42250
$ws.Cells.Item(87, 12).Value = 42250
$ws.Cells.Item(87, 14).Value = -44496

$ws.Cells.Item(88, 8).Value = 40990
$ws.Cells.Item(88, 10).Value = 40990
$ws.Cells.Item(88, 12).Value = 40990
$ws.Cells.Item(88, 14).Value = -41846

$ws.Cells.Item(90, 8).Value = 42250
$ws.Cells.Item(90, 10).Value = 42250
$ws.Cells.Item(90, 12).Value = 126750
$ws.Cells.Item(90, 14).Value = -137982

$ws.Cells.Item(91, 8).Value = 40990
$ws.Cells.Item(91, 10).Value = 40990
$ws.Cells.Item(91, 12).Value = 40990
$ws.Cells.Item(91, 14).Value = -43954

$ws.Cells.Item(126, 8).Value = 2589.5789
$ws.Cells.Item(126, 9).Value = 2587.4375
$ws.Cells.Item(126, 10).Value = 2601
$ws.Cells.Item(126, 11).Value = 7762.3125
$ws.Cells.Item(126, 12).Value = 7803
$ws.Cells.Item(126, 13).Value = -5292.3125
$ws.Cells.Item(126, 14).Value = -12743

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 507.65
$ws.Cells.Item(107, 9).Value = 508.05264
$ws.Cells.Item(107, 11).Value = 1524.15792
$ws.Cells.Item(107, 13).Value = 395.8420799999999
